$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 42, pushing the existing rows 42-45 down to 43-46
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly entry
$ws.Range("A42").Value = 11
$ws.Range("B42").Value = "Vega Monumental Concepción"
$ws.Range("C42").Value = "Bíobío"
$ws.Range("D42").Value = 44476
$ws.Range("E42").Value = 8
$ws.Range("F42").Value = 100112012
$ws.Range("G42").Value = "Espinaca"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 80
$ws.Range("K42").Value = 5000
$ws.Range("L42").Value = 5500
$ws.Range("M42").Value = 5312
$ws.Range("N42").Value = "$/cuna 10 kilos"
$ws.Range("O42").Value = "Región Metropolitana"
$ws.Range("P42").Value = 531
$ws.Range("Q42").Value = 10
$ws.Range("R42").Value = "Hortaliza"
